$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (changed) date column C for rows 2-5 from 45204 (2023-10-05) to 45207 (2023-10-08)
$ws.Range("C2:C5").Value = 45207
